$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_1_9_0"
$ws.Range("B2").Value = 0.9903547028684991
$ws.Range("C2").Value = 0.9996339600147776
$ws.Range("D2").Value = 0.9982938703912183
$ws.Range("E2").Value = 0.9988000297943297
$ws.Range("F2").Value = 1.163225705759224
$ws.Range("G2").Value = 0.02797977922273666
$ws.Range("H2").Value = 0.2419368743981863
$ws.Range("I2").Value = 0.1286652799108416
$ws.Range("A3").Value = "model_1_9_1"
$ws.Range("B3").Value = 0.9919720579355564
$ws.Range("C3").Value = 0.9988141059215304
$ws.Range("D3").Value = 0.9978736892578257
$ws.Range("E3").Value = 0.9982290973833838
$ws.Range("F3").Value = 0.9681722031359948
$ws.Range("G3").Value = 0.09064871554117977
$ws.Range("H3").Value = 0.3015204544327065
$ws.Range("I3").Value = 0.1898827819099711
$ws.Range("A4").Value = "model_1_9_2"
$ws.Range("B4").Value = 0.9932167450461512
$ws.Range("C4").Value = 0.99768407283584
$ws.Range("D4").Value = 0.9974262682729748
$ws.Range("E4").Value = 0.9975241565826702
$ws.Range("F4").Value = 0.8180625670168155
$ws.Range("G4").Value = 0.1770274652091698
$ws.Range("H4").Value = 0.3649667682753586
$ws.Range("I4").Value = 0.26546916315159
$ws.Range("A5").Value = "model_1_9_3"
$ws.Range("B5").Value = 0.9941687795631367
$ws.Range("C5").Value = 0.9963742167005276
$ws.Range("D5").Value = 0.9969688114111923
$ws.Range("E5").Value = 0.9967450908459101
$ws.Range("F5").Value = 0.7032469208185538
$ws.Range("G5").Value = 0.2771517329372236
$ws.Range("H5").Value = 0.4298362147359261
$ws.Range("I5").Value = 0.349003496433797
$ws.Range("A6").Value = "model_1_9_4"
$ws.Range("B6").Value = 0.9948911715090366
$ws.Range("C6").Value = 0.9949778797509877
$ws.Range("D6").Value = 0.9965117630917886
$ws.Range("E6").Value = 0.9959336517349833
$ws.Range("F6").Value = 0.6161262370648239
$ws.Range("G6").Value = 0.383886519151714
$ws.Range("H6").Value = 0.4946477280443722
$ws.Range("I6").Value = 0.4360090236082731
$ws.Range("A7").Value = "model_1_9_5"
$ws.Range("B7").Value = 0.9954335821709813
$ws.Range("C7").Value = 0.993560017430285
$ws.Range("D7").Value = 0.9960615108316353
$ws.Range("E7").Value = 0.9951182919344899
$ws.Range("F7").Value = 0.5507113497420221
$ws.Range("G7").Value = 0.492266686081811
$ws.Range("H7").Value = 0.5584955294959721
$ws.Range("I7").Value = 0.5234349417375695
$ws.Range("A8").Value = "model_1_9_6"
$ws.Range("B8").Value = 0.9958351409715774
$ws.Range("C8").Value = 0.9921644434443859
$ws.Range("D8").Value = 0.9956219874080805
$ws.Range("E8").Value = 0.9943180530710185
$ws.Range("F8").Value = 0.502283238833783
$ws.Range("G8").Value = 0.5989431520168598
$ws.Range("H8").Value = 0.6208219335231411
$ws.Range("I8").Value = 0.609239536616292
$ws.Range("A9").Value = "model_1_9_7"
$ws.Range("B9").Value = 0.9961267040587813
$ws.Range("C9").Value = 0.9908196386542313
$ws.Range("D9").Value = 0.9951957497179718
$ws.Range("E9").Value = 0.9935452102359048
$ws.Range("F9").Value = 0.4671206437097407
$ws.Range("G9").Value = 0.7017388646309803
$ws.Range("H9").Value = 0.6812643606194445
$ws.Range("I9").Value = 0.6921066271799684
$ws.Range("A10").Value = "model_1_9_24"
$ws.Range("B10").Value = 0.9961952652142357
$ws.Range("C10").Value = 0.9786340144819659
$ws.Range("D10").Value = 0.9906588490336753
$ws.Range("E10").Value = 0.9861226241666595
$ws.Range("F10").Value = 0.4588521479491812
$ws.Range("G10").Value = 1.633197415051405
$ws.Range("H10").Value = 1.324617342341333
$ws.Range("I10").Value = 1.487983982925016
$ws.Range("A11").Value = "model_1_9_23"
$ws.Range("B11").Value = 0.9962351090192327
$ws.Range("C11").Value = 0.9789377885706064
$ws.Range("D11").Value = 0.9907945414635797
$ws.Range("E11").Value = 0.9863217278238479
$ws.Range("F11").Value = 0.4540469732038327
$ws.Range("G11").Value = 1.609977186997397
$ws.Range("H11").Value = 1.305375543710334
$ws.Range("I11").Value = 1.466635346381918
$ws.Range("A12").Value = "model_1_9_22"
$ws.Range("B12").Value = 0.9962771083441549
$ws.Range("C12").Value = 0.9792737109163676
$ws.Range("D12").Value = 0.9909430067586844
$ws.Range("E12").Value = 0.9865408916626116
$ws.Range("F12").Value = 0.4489818421137216
$ws.Range("G12").Value = 1.584299574032052
$ws.Range("H12").Value = 1.284322495179111
$ws.Range("I12").Value = 1.443135782369725
$ws.Range("A13").Value = "model_1_9_21"
$ws.Range("B13").Value = 0.9963210231215068
$ws.Range("C13").Value = 0.9796453438414657
$ws.Range("D13").Value = 0.9911051270165158
$ws.Range("E13").Value = 0.9867820584559616
$ws.Range("F13").Value = 0.4436857068903098
$ws.Range("G13").Value = 1.55589227533263
$ws.Range("H13").Value = 1.261333111339535
$ws.Range("I13").Value = 1.417276979521999
$ws.Range("A14").Value = "model_1_9_8"
$ws.Range("B14").Value = 0.9963326694078716
$ws.Range("C14").Value = 0.9895428067642321
$ws.Range("D14").Value = 0.9947846265944545
$ws.Range("E14").Value = 0.9928074815903376
$ws.Range("F14").Value = 0.4422811612872651
$ws.Range("G14").Value = 0.7993387876694829
$ws.Range("H14").Value = 0.7395634739954697
$ws.Range("I14").Value = 0.7712086434063763
$ws.Range("A15").Value = "model_1_9_20"
$ws.Range("B15").Value = 0.996366355803611
$ws.Range("C15").Value = 0.9800557155858439
$ws.Range("D15").Value = 0.9912817758864158
$ws.Range("E15").Value = 0.98704689858419
$ws.Range("F15").Value = 0.4382185719316231
$ws.Range("G15").Value = 1.524523814862468
$ws.Range("H15").Value = 1.236283504774121
$ws.Range("I15").Value = 1.388879833435276
$ws.Range("A16").Value = "model_1_9_19"
$ws.Range("B16").Value = 0.9964125216072826
$ws.Range("C16").Value = 0.9805088744123874
$ws.Range("D16").Value = 0.9914739761735065
$ws.Range("E16").Value = 0.9873375285039129
$ws.Range("F16").Value = 0.4326509622638599
$ws.Range("G16").Value = 1.489884746915246
$ws.Range("H16").Value = 1.209028637102971
$ws.Range("I16").Value = 1.357717409739331
$ws.Range("A17").Value = "model_1_9_18"
$ws.Range("B17").Value = 0.9964585935345996
$ws.Range("C17").Value = 0.9810085958492913
$ws.Range("D17").Value = 0.9916824939925423
$ws.Range("E17").Value = 0.9876558948677517
$ws.Range("F17").Value = 0.4270946741124036
$ws.Range("G17").Value = 1.451686473387991
$ws.Range("H17").Value = 1.1794598697982
$ws.Range("I17").Value = 1.323580981081383
$ws.Range("A18").Value = "model_1_9_9"
$ws.Range("B18").Value = 0.996472329500017
$ws.Range("C18").Value = 0.9883436596659108
$ws.Range("D18").Value = 0.9943900934571762
$ws.Range("E18").Value = 0.992109369251766
$ws.Range("F18").Value = 0.4254381125934551
$ws.Range("G18").Value = 0.8910005525617104
$ws.Range("H18").Value = 0.795510052490047
$ws.Range("I18").Value = 0.8460628514751138
$ws.Range("A19").Value = "model_1_9_17"
$ws.Range("B19").Value = 0.9965034099978137
$ws.Range("C19").Value = 0.9815595000022879
$ws.Range("D19").Value = 0.9919080855763915
$ws.Range("E19").Value = 0.9880042541404788
$ws.Range("F19").Value = 0.4216897953055499
$ws.Range("G19").Value = 1.409575837402792
$ws.Range("H19").Value = 1.147469965627902
$ws.Range("I19").Value = 1.286228600894641
$ws.Range("A20").Value = "model_1_9_16"
$ws.Range("B20").Value = 0.9965453330389757
$ws.Range("C20").Value = 0.9821657562135409
$ws.Range("D20").Value = 0.9921517027134228
$ws.Range("E20").Value = 0.9883846662220286
$ws.Range("F20").Value = 0.4166338640596372
$ws.Range("G20").Value = 1.36323413805822
$ws.Range("H20").Value = 1.11292395670818
$ws.Range("I20").Value = 1.245439399027143
$ws.Range("A21").Value = "model_1_9_10"
$ws.Range("B21").Value = 0.9965609383744592
$ws.Range("C21").Value = 0.9872266878804726
$ws.Range("D21").Value = 0.994013226115905
$ws.Range("E21").Value = 0.9914532347411728
$ws.Range("F21").Value = 0.4147518559541362
$ws.Range("G21").Value = 0.9763809077586699
$ws.Range("H21").Value = 0.8489515414253676
$ws.Range("I21").Value = 0.9164160403007306
$ws.Range("A22").Value = "model_1_9_15"
$ws.Range("B22").Value = 0.9965822321818147
$ws.Range("C22").Value = 0.9828321267343628
$ws.Range("D22").Value = 0.9924138039521166
$ws.Range("E22").Value = 0.9887992829735787
$ws.Range("F22").Value = 0.4121838164472477
$ws.Range("G22").Value = 1.312297352991397
$ws.Range("H22").Value = 1.075756818796089
$ws.Range("I22").Value = 1.200982644899562
$ws.Range("A23").Value = "model_1_9_11"
$ws.Range("B23").Value = 0.9966105687573016
$ws.Range("C23").Value = 0.9861922706547841
$ws.Range("D23").Value = 0.9936547627691673
$ws.Range("E23").Value = 0.9908397863166917
$ws.Range("F23").Value = 0.4087664170068745
$ws.Range("G23").Value = 1.055450863958598
$ws.Range("H23").Value = 0.8997832609205628
$ws.Range("I23").Value = 0.982192267805172
$ws.Range("A24").Value = "model_1_9_14"
$ws.Range("B24").Value = 0.9966113459779248
$ws.Range("C24").Value = 0.9835633744210227
$ws.Range("D24").Value = 0.9926949372241863
$ws.Range("E24").Value = 0.9892502375597699
$ws.Range("F24").Value = 0.4086726839683044
$ws.Range("G24").Value = 1.256401413597107
$ws.Range("H24").Value = 1.035890852703145
$ws.Range("I24").Value = 1.152629612644932
$ws.Range("A25").Value = "model_1_9_13"
$ws.Range("B25").Value = 0.996629063346517
$ws.Range("C25").Value = 0.9843642465788082
$ws.Range("D25").Value = 0.9929954555362295
$ws.Range("E25").Value = 0.9897395128680122
$ws.Range("F25").Value = 0.4065359640410806
$ws.Range("G25").Value = 1.195183439973666
$ws.Range("H25").Value = 0.9932760005014643
$ws.Range("I25").Value = 1.100167689681354
$ws.Range("A26").Value = "model_1_9_12"
$ws.Range("B26").Value = 0.9966307743828222
$ws.Range("C26").Value = 0.9852391845488315
$ws.Range("D26").Value = 0.9933154542406443
$ws.Range("E26").Value = 0.9902688843259325
$ws.Range("F26").Value = 0.406329612553252
$ws.Range("G26").Value = 1.128303939855769
$ws.Range("H26").Value = 0.9478987407908841
$ws.Range("I26").Value = 1.043406508038456